# Refresh the crypto price/volume snapshot (GitHub Actions bot update).
#
# Notes:
#  - Column D ("Price") holds numeric-looking text (e.g. "63.690.23",
#    "570.62", "1.00") that must stay plain text, exactly as authored in the
#    source feed. Assigning a bare numeric-looking string to Range.Value
#    lets Excel auto-convert it to a real number, so every Column D write is
#    prefixed with a literal leading apostrophe (the same trick Excel's UI
#    uses for "force text") to keep it a string.
#  - Column E ("Volume(1h)") values already contain a '%' plus padding
#    spaces, so Excel always stores them as text - no extra handling needed.
#  - Rows 26 and 27 swap their Coin/Link/Price/Volume contents
#    (Binance-PegBSC-USD <-> Fetch.AI) while keeping the same rank index in
#    column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'63.690.23"
$ws.Range("E2").Value  = "  +6.52%  "

$ws.Range("D3").Value  = "'2.474.51"
$ws.Range("E3").Value  = "  +7.62%  "

$ws.Range("E4").Value  = "  -0.06%  "

$ws.Range("D5").Value  = "'570.62"
$ws.Range("E5").Value  = "  +5.61%  "

$ws.Range("D6").Value  = "'143.72"
$ws.Range("E6").Value  = "  +11.89%  "

$ws.Range("D7").Value  = "'1.00"
$ws.Range("E7").Value  = "  -0.04%  "

$ws.Range("D8").Value  = "'0.593"
$ws.Range("E8").Value  = "  +4.50%  "

$ws.Range("D9").Value  = "'2.473.47"
$ws.Range("E9").Value  = "  +7.65%  "

$ws.Range("E10").Value = "  +6.42%  "

$ws.Range("D11").Value = "'5.78"
$ws.Range("E11").Value = "  +4.97%  "

$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("D13").Value = "'0.356"
$ws.Range("E13").Value = "  +7.97%  "

$ws.Range("D14").Value = "'26.52"
$ws.Range("E14").Value = "  +14.98%  "

$ws.Range("D15").Value = "'2.918.99"
$ws.Range("E15").Value = "  +7.52%  "

$ws.Range("D16").Value = "'63.541.12"
$ws.Range("E16").Value = "  +6.34%  "

$ws.Range("E17").Value = "  +10.02%  "

$ws.Range("D18").Value = "'2.473.00"
$ws.Range("E18").Value = "  +7.30%  "

$ws.Range("D19").Value = "'11.35"
$ws.Range("E19").Value = "  +9.21%  "

$ws.Range("D20").Value = "'344.90"
$ws.Range("E20").Value = "  +11.30%  "

$ws.Range("D21").Value = "'4.34"
$ws.Range("E21").Value = "  +8.74%  "

$ws.Range("E22").Value = "  +6.68%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "'65.93"
$ws.Range("E24").Value = "  +4.06%  "

$ws.Range("D25").Value = "'0.177"
$ws.Range("E25").Value = "  +4.97%  "

# Row 26 becomes what Row 27 used to be (Fetch.AI)
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "'1.54"
$ws.Range("E26").Value = "  +15.40%  "

# Row 27 becomes what Row 26 used to be (Binance-PegBSC-USD)
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").Value = "'8.26"
$ws.Range("E28").Value = "  +7.04%  "

$ws.Range("D29").Value = "'1.34"
$ws.Range("E29").Value = "  +14.04%  "

$ws.Range("D30").Value = "'0.0" + [char]0x2083 + "0822"
$ws.Range("E30").Value = "  +15.77%  "

$ws.Range("E31").Value = "  +19.23%  "

$ws.Range("E32").Value = "  +8.99%  "

$ws.Range("D33").Value = "'175.87"
$ws.Range("E33").Value = "  +2.61%  "

$ws.Range("D34").Value = "'1.53"
$ws.Range("E34").Value = "  +13.12%  "

$ws.Range("D35").Value = "'0.402"
$ws.Range("E35").Value = "  +6.33%  "

$ws.Range("D36").Value = "'19.10"
$ws.Range("E36").Value = "  +7.89%  "

$ws.Range("D37").Value = "'373.68"
$ws.Range("E37").Value = "  +20.90%  "

$ws.Range("E38").Value = "  +12.25%  "

$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "'1.74"
$ws.Range("E40").Value = "  +16.33%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").Value = "'40.37"
$ws.Range("E42").Value = "  +6.58%  "

$ws.Range("D43").Value = "'152.15"
$ws.Range("E43").Value = "  +11.89%  "

$ws.Range("E44").Value = "  +10.57%  "

$ws.Range("D45").Value = "'20.84"
$ws.Range("E45").Value = "  +12.96%  "

$ws.Range("D46").Value = "'0.600"
$ws.Range("E46").Value = "  +6.68%  "

$ws.Range("E47").Value = "  +3.75%  "

$ws.Range("E48").Value = "  +8.25%  "

$ws.Range("D49").Value = "'0.0" + [char]0x2086 + "0236"
$ws.Range("E49").Value = "  +7.45%  "

$ws.Range("D50").Value = "'0.0228"
$ws.Range("E50").Value = "  +8.06%  "

$ws.Range("D51").Value = "'18.31"
$ws.Range("E51").Value = "  +10.37%  "
